# This workbook is a weekly re-shuffle of the daily price rows: the values
# in columns D (Fecha), J (Volumen), K/L/M (Precios), O (Origen) and
# P (Precio $/Kg) for each data row (2..19) are replaced by the values that
# used to live in a *different* row, according to the permutation below.
# Columns A, B, C, E, F, G, H, I, N, Q, R are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row -> source row (values are taken from the ORIGINAL
# workbook at the source row and written into the destination row).
$rowMap = @{
    2  = 9
    3  = 16
    4  = 5
    5  = 19
    6  = 7
    7  = 15
    8  = 14
    9  = 2
    10 = 13
    11 = 6
    12 = 11
    13 = 8
    14 = 3
    15 = 17
    16 = 12
    17 = 4
    18 = 10
    19 = 18
}

$cols = @("D", "J", "K", "L", "M", "O", "P")

# First, snapshot the original values for every relevant cell, since rows
# overwrite each other as we iterate (it's a permutation/cycle).
# NOTE: use .Value2 (not .Value) when reading, since .Value alone does not
# resolve to the underlying scalar in this runtime.
$original = @{}
foreach ($row in 2..19) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $original[$addr] = $ws.Range($addr).Value2
    }
}

# Now write the shuffled values into each destination row.
foreach ($destRow in 2..19) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $srcAddr = "$col$srcRow"
        $destAddr = "$col$destRow"
        $ws.Range($destAddr).Value = $original[$srcAddr]
    }
}
